$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# New "Programs/Studies/Cases/Samples/Case Files/Study Files" stat query,
# replacing the old per-tab "number_of_*" style stat query that used to be
# shared (as the StatQuery column) across all three tab rows.
$programsQuery = "MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)`nOPTIONAL MATCH (samp:sample)-->(c)`nOPTIONAL MATCH (diag:diagnosis)-->(c)`nOPTIONAL MATCH (f:file)-[*]->(c)`nOPTIONAL MATCH (sf:file)-->(s)`nWITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p`nWHERE diag.stage_of_disease IN ['IVb']`nRETURN  `n    count(distinct p) AS Programs,`n    count(distinct s) AS Studies,`n    count(distinct c) AS Cases,`n    count(distinct samp) AS Samples,`n    count(distinct f) AS ``Case Files``,`n    count(distinct sf) AS ``Study Files``"

# Cases-tab query (row 2, column B) gains a trailing "Cohort" projection.
$caseQuery = "MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)`nMATCH (c)<--(diag:diagnosis)`nOPTIONAL MATCH (samp:sample)-->(c)`nOPTIONAL MATCH (co:cohort)<-[*]-(c)`nWITH DISTINCT c, s, demo, diag, co`nWHERE diag.stage_of_disease IN ['IVb']`nRETURN  coalesce(c.case_id, '') AS ``Case ID`` ,`n        coalesce(s.clinical_study_designation, '') AS ``Study Code`` ,`n        coalesce(s.clinical_study_type, '') AS  ``Study Type``,`n        coalesce(demo.breed, '') AS Breed ,`n        coalesce(diag.disease_term, '') AS Diagnosis ,`n        coalesce(diag.stage_of_disease, '') AS ``Stage of Disease`` ,`n        coalesce(demo.patient_age_at_enrollment, '') AS Age ,`n        coalesce(demo.sex, '') AS Sex ,`n        coalesce(demo.neutered_indicator, '') AS ``Neutered Status``,`n        coalesce(demo.weight, '') AS ``Weight (kg)``,`n        coalesce(diag.best_response, '') AS ``Response to Treatment``,`ncoalesce(co.cohort_description, '') AS ``Cohort``"

# Samples-tab query (row 3, column B) is unchanged.
$sampleQuery = "MATCH (s:study)<-[*]-(c:case)<--(demo:demographic), (samp:sample)-->(c)<--(diag:diagnosis) `nWHERE diag.stage_of_disease IN ['IVb']`nWITH DISTINCT samp AS samp, c, demo, diag`nRETURN  coalesce(samp.sample_id, '') AS ``Sample ID``, `n        coalesce(c.case_id, '') AS ``Case ID``, `n        coalesce(demo.breed,'') AS Breed , `n        coalesce(diag.disease_term,'') AS Diagnosis , `n        coalesce(samp.sample_site, '') AS ``Sample Site``,`n        coalesce(samp.summarized_sample_type, '') AS ``Sample Type``,`n        coalesce(samp.specific_sample_pathology, '') AS ``Pathology/Morphology``,`n        coalesce(samp.tumor_grade, '') AS ``Tumor Grade``,`n        coalesce(samp.sample_chronology, '') AS ``Sample Chronology``,`n        coalesce(samp.percentage_tumor, '') AS ``Percentage Tumor``,`n        coalesce(samp.necropsy_sample, '') AS ``Necropsy Sample``,`n        coalesce(samp.sample_preservation, '') AS ``Sample Preservation``"

# Files-tab query (row 4, column B) drops the trailing "Study Code" projection.
$fileQuery = "`nMATCH (f:file)-->(parent)`nWITH DISTINCT f, parent`nMATCH (f)-[*]->(c:case)<--(demo:demographic)`n MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)`nWHERE diag.stage_of_disease IN ['IVb']`nWITH DISTINCT f, parent, c, demo, diag, s`nRETURN coalesce(f.file_name, '') AS ``File Name``, `n        coalesce(labels(parent)[0], '') AS ``Association``,`n        coalesce(f.file_description, '') AS ``Description``,`n        coalesce(f.file_format, '') AS ``Format``,`n        coalesce(f.file_size, '') AS ``Size``,`n        coalesce(c.case_id, '') AS ``Case ID``, `n        coalesce(diag.disease_term,'') AS Diagnosis "

$ws.Range("B2").Value = $caseQuery
$ws.Range("C2").Value = $programsQuery

$ws.Range("B3").Value = $sampleQuery
$ws.Range("C3").Value = $programsQuery

$ws.Range("B4").Value = $fileQuery
$ws.Range("C4").Value = $programsQuery

# Update active selection to C2, matching the new sheetView.
$ws.Range("C2").Select()

$wb.Save()
